$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Text change: "Ready for handoff" -> "In Translation" ---
# Every cell that held the old status text must be updated so the shared
# string is fully replaced (no remaining reference to the old text).
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width change: 17.2159881591797 -> 13.4101845877511 ---
# The ColumnWidth setter here quantizes to the nearest 1/6 of a character
# (internal storage granularity), so feed it the pre-quantized character
# width (12.5) that lands closest to the target OOXML width after the
# engine's own "+5/6" padding is re-applied.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
